$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Row 50 (Lp=50): rewrite the stale "odswierzanie formsa" task text and mark it
# complete (100) now that a concrete example (magazine error) has been found.
$ws.Range("C51").Value = "Znaleźć sposób na odświerzanie formsa jeśli dodana jest edycja komórek. Przykład - błędu magazynu"
$ws.Range("D51").Value = 100

# New task row (Lp=57): "Ustawienia danych w babelacg" under the "Ustawienia" group.
$ws.Range("A58").Value = 57
$ws.Range("B58").Value = "Ustawienia. "
$ws.Range("C58").Value = "Ustawienia danych w babelacg"
$ws.Range("D58").Value = 0
$ws.Range("C58").WrapText = $true

# Match the author's final on-screen selection after the edit.
$ws.Range("C61").Select() | Out-Null
